# Weekly update: add two new price rows (week of 2023-10-13) for
# "Repollo" at Feria Lagunitas de Puerto Montt, pushing the existing
# historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 756, shifting every
# row from 756 downward (including the last data row, 835) down by two.
$ws.Range("A756:R757").Insert()

# New row 756: Copenhague / Primera
$ws.Range("A756").Value = 4
$ws.Range("B756").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C756").Value = "Los Lagos"
$ws.Range("D756").Value = 45212
$ws.Range("E756").Value = 10
$ws.Range("F756").Value = 100112006
$ws.Range("G756").Value = "Repollo"
$ws.Range("H756").Value = "Copenhague"
$ws.Range("I756").Value = "Primera"
$ws.Range("J756").Value = 500
$ws.Range("K756").Value = 1800
$ws.Range("L756").Value = 1800
$ws.Range("M756").Value = 1800
$ws.Range("N756").Value = '$/unidad'
$ws.Range("O756").Value = "Región Metropolitana"
$ws.Range("P756").Value = 1800
$ws.Range("Q756").Value = 1
$ws.Range("R756").Value = "Hortaliza"

# New row 757: Crespo record / Primera
$ws.Range("A757").Value = 4
$ws.Range("B757").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C757").Value = "Los Lagos"
$ws.Range("D757").Value = 45212
$ws.Range("E757").Value = 10
$ws.Range("F757").Value = 100112006
$ws.Range("G757").Value = "Repollo"
$ws.Range("H757").Value = "Crespo record"
$ws.Range("I757").Value = "Primera"
$ws.Range("J757").Value = 1000
$ws.Range("K757").Value = 1500
$ws.Range("L757").Value = 1500
$ws.Range("M757").Value = 1500
$ws.Range("N757").Value = '$/unidad'
$ws.Range("O757").Value = "Región Metropolitana"
$ws.Range("P757").Value = 1500
$ws.Range("Q757").Value = 1
$ws.Range("R757").Value = "Hortaliza"
